$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "45.116.34"
$ws.Range("E2").Value = "  +0.82%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.272.33"
$ws.Range("E4").Value = "  -0.87%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "302.00"
$ws.Range("E5").Value = "  -1.83%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "94.90"
$ws.Range("E6").Value = "  -1.49%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.566"
$ws.Range("E7").Value = "  -1.24%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.998"
$ws.Range("E8").Value = "  -0.58%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.511"
$ws.Range("E9").Value = "  -2.62%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.45"
$ws.Range("E10").Value = "  -2.48%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0790"
$ws.Range("E11").Value = "  -2.20%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.23"
$ws.Range("E12").Value = "  -0.75%  "
$ws.Range("E13").Value = "  -0.61%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.615.42"
$ws.Range("E14").Value = "  +1.01%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.269.42"
$ws.Range("E15").Value = "  +0.78%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.79"
$ws.Range("E16").Value = "  +0.92%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.799"
$ws.Range("E17").Value = "  -5.05%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "44.887.41"
$ws.Range("E18").Value = "  +0.87%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.99"
$ws.Range("E19").Value = "  +7.24%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0926"
$ws.Range("E20").Value = "  -3.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.11"
$ws.Range("E21").Value = "  -3.65%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "65.40"
$ws.Range("E22").Value = "  -0.48%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "239.61"
$ws.Range("E23").Value = "  +0.27%  "
$ws.Range("E24").Value = "  -3.16%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.998"
$ws.Range("E25").Value = "  -0.59%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.94"
$ws.Range("E26").Value = "  -4.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "42.24"
$ws.Range("E27").Value = "  +11.98%  "
$ws.Range("E28").Value = "  +1.26%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.56"
$ws.Range("E29").Value = "  -3.22%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.60"
$ws.Range("E30").Value = "  -2.22%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "152.39"
$ws.Range("E31").Value = "  -0.27%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.62"
$ws.Range("E32").Value = "  -7.19%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0788"
$ws.Range("E33").Value = "  -1.80%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.57"
$ws.Range("E34").Value = "  -2.78%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.96"
$ws.Range("E35").Value = "  -4.01%  "
$ws.Range("E36").Value = "  -1.66%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.106"
$ws.Range("E37").Value = "  -3.95%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.76"
$ws.Range("E38").Value = "  -5.82%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.87"
$ws.Range("E39").Value = "  +2.18%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0306"
$ws.Range("E40").Value = "  +1.36%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.28"
$ws.Range("E41").Value = "  -4.43%  "
$ws.Range("E42").Value = "  -9.12%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("E43").Value = "  -0.77%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.97"
$ws.Range("E44").Value = "  +13.17%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.779.22"
$ws.Range("E45").Value = "  -2.62%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.195"
$ws.Range("E46").Value = "  +0.37%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "70.47"
$ws.Range("E47").Value = "  -0.93%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "75.84"
$ws.Range("E48").Value = "  -5.05%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "97.12"
$ws.Range("E49").Value = "  -2.50%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.87"
$ws.Range("E50").Value = "  -2.66%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "53.31"
$ws.Range("E51").Value = "  -3.11%  "
